# "Fruta / hortaliza, semanal" update for the Zanahoria (Macroferia Regional
# de Talca) sheet: a new weekly record is inserted at row 122, pushing the
# existing records (previously rows 122-183) down to rows 123-184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 122 - this shifts rows
# 122..183 down to 123..184 and grows the sheet dimension to A1:R184.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new weekly observation.
$ws.Range("A122").Value = 5
$ws.Range("B122").Value = "Macroferia Regional de Talca"
$ws.Range("C122").Value = "Maule"
$ws.Range("D122").Value = 44452
$ws.Range("E122").Value = 7
$ws.Range("F122").Value = 100114013
$ws.Range("G122").Value = "Zanahoria"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 5500
$ws.Range("L122").Value = 5500
$ws.Range("M122").Value = 5500
$ws.Range("N122").Value = "$/saco 20 kilos"
$ws.Range("O122").Value = "Región de Ñuble"
$ws.Range("P122").Value = 275
$ws.Range("Q122").Value = 20
$ws.Range("R122").Value = "Hortaliza"
